# fix: Q13 and Q15 marked as 基本法 while it should be 憲法
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Q13 is row 14 (ID=13), Q15 is row 16 (ID=15): QuestionType column E
$ws.Range("E14").Value = "憲法"
$ws.Range("E16").Value = "憲法"

# Update the active cell selection to match the saved view state
$ws.Range("C25").Select()
